$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F = 想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 304
$ws1.Range("F3").Value = 1213
$ws1.Range("F4").Value = 16900
$ws1.Range("F8").Value = 7
$ws1.Range("F12").Value = 11719
$ws1.Range("F14").Value = 1393
$ws1.Range("F15").Value = 4656
$ws1.Range("F16").Value = 466
$ws1.Range("F17").Value = 5
$ws1.Range("F18").Value = 407
$ws1.Range("F20").Value = 900
$ws1.Range("F21").Value = 341
$ws1.Range("F24").Value = 5214

# Sheet "全部类型" updates (column F = 想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 304
$ws4.Range("F4").Value = 1213
$ws4.Range("F5").Value = 16900
$ws4.Range("F9").Value = 7
$ws4.Range("F15").Value = 11719
$ws4.Range("F17").Value = 1393
$ws4.Range("F18").Value = 4656
$ws4.Range("F19").Value = 466
$ws4.Range("F20").Value = 5
$ws4.Range("F21").Value = 407
$ws4.Range("F23").Value = 900
$ws4.Range("F24").Value = 341
$ws4.Range("F27").Value = 5214
